$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47; this shifts the existing rows 47-50
# down to 48-51, preserving their contents and formatting.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the latest weekly price entry.
$ws.Range("A47").Value2 = 1
$ws.Range("B47").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C47").Value2 = "Arica y Parinacota"
$ws.Range("D47").Value2 = 44931
$ws.Range("E47").Value2 = 15
$ws.Range("F47").Value2 = 100112045
$ws.Range("G47").Value2 = "Zapallo"
$ws.Range("H47").Value2 = "Camote"
$ws.Range("I47").Value2 = "1a nueva(o)"
$ws.Range("J47").Value2 = 1100
$ws.Range("K47").Value2 = 750
$ws.Range("L47").Value2 = 760
$ws.Range("M47").Value2 = 755
$ws.Range("N47").Value2 = "$/kilo (volumen en unidades)"
$ws.Range("O47").Value2 = "Perú"
$ws.Range("P47").Value2 = 755
$ws.Range("Q47").Value2 = 1
$ws.Range("R47").Value2 = "Hortaliza"
